$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9: date, time, and text for "everything else" work
# Reuse existing number formats from the prior row (B8 = date, C8 = time)
$ws.Range("B8:C8").Copy()
$ws.Range("B9:C9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B9").Value = 42160
$ws.Range("C9").Value = 0.33333333333333331

$ws.Range("E9").Value = "everything else"

$ws.Range("B10").Select()
